# "Refined metadata to be additional tab"
#
# The former single-sheet workbook ("data") gets a second tab, "metadata",
# appended after it, carrying the panel-query bookkeeping fields
# (data_name / data_id / data_version / data_version_created /
#  panel_query_time / panel_get_request). The "data" sheet's per-row
# time_taken column (F) is also refreshed to the timestamps captured by
# the same (re-run) scrape.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- refresh the "time_taken" column on the "data" sheet --------------------
$newTimes = @(
    "2021-10-05 14:21:09.273897",
    "2021-10-05 14:21:09.273906",
    "2021-10-05 14:21:09.273909",
    "2021-10-05 14:21:09.273912",
    "2021-10-05 14:21:09.273915",
    "2021-10-05 14:21:09.273918",
    "2021-10-05 14:21:09.273920",
    "2021-10-05 14:21:09.273923",
    "2021-10-05 14:21:09.273926",
    "2021-10-05 14:21:09.273929",
    "2021-10-05 14:21:09.273932",
    "2021-10-05 14:21:09.273934",
    "2021-10-05 14:21:09.273937",
    "2021-10-05 14:21:09.273939",
    "2021-10-05 14:21:09.273942",
    "2021-10-05 14:21:09.273944",
    "2021-10-05 14:21:09.273947",
    "2021-10-05 14:21:09.273950",
    "2021-10-05 14:21:09.273952",
    "2021-10-05 14:21:09.273955"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = 2 + $i
    $dataSheet.Range("F$row").Value = $newTimes[$i]
}

# --- add the new "metadata" tab right after "data" --------------------------
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# header row (bold, centered, thin-bordered) mirrors the style used for the
# "data" sheet's own header row / index column
$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$indexCell = $ws.Range("A2")
$indexCell.Font.Bold = $true
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
$indexCell.Borders.LineStyle = 1
$indexCell.Value = 0

$ws.Range("B2").Value = "Inherited phaeochromocytoma and paraganglioma"
$ws.Range("C2").Value = 97

# data_version must stay textual ("1.9"), not become the number 1.9
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.9"

$ws.Range("E2").Value = "2021-08-02T08:09:28.774152Z"
$ws.Range("F2").Value = "2021-10-05 14:21:09.270727"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/97/?format=json"

$ws.Range("A1").Select()

# keep "data" as the active tab, same as before the edit
$dataSheet.Activate()
